# Applies the F-column (想去人数 / want-to-go count) updates and the
# two cover-image URL swaps (I40 on 展览, I44 on 全部类型) described by
# the commit's OOXML diff (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook
$mismatches = 0
$applied = 0

$ws1 = $wb.Worksheets.Item("展览")
$cur = $ws1.Range("F3").Value()
if ($cur -ne 113) { $mismatches++ }
$ws1.Range("F3").Value = 114
$applied++
$cur = $ws1.Range("F4").Value()
if ($cur -ne 555) { $mismatches++ }
$ws1.Range("F4").Value = 557
$applied++
$cur = $ws1.Range("F6").Value()
if ($cur -ne 18) { $mismatches++ }
$ws1.Range("F6").Value = 20
$applied++
$cur = $ws1.Range("F7").Value()
if ($cur -ne 1907) { $mismatches++ }
$ws1.Range("F7").Value = 1910
$applied++
$cur = $ws1.Range("F8").Value()
if ($cur -ne 5051) { $mismatches++ }
$ws1.Range("F8").Value = 5068
$applied++
$cur = $ws1.Range("F9").Value()
if ($cur -ne 1412) { $mismatches++ }
$ws1.Range("F9").Value = 1419
$applied++
$cur = $ws1.Range("F10").Value()
if ($cur -ne 148) { $mismatches++ }
$ws1.Range("F10").Value = 149
$applied++
$cur = $ws1.Range("F11").Value()
if ($cur -ne 3005) { $mismatches++ }
$ws1.Range("F11").Value = 3010
$applied++
$cur = $ws1.Range("F13").Value()
if ($cur -ne 29) { $mismatches++ }
$ws1.Range("F13").Value = 31
$applied++
$cur = $ws1.Range("F14").Value()
if ($cur -ne 1250) { $mismatches++ }
$ws1.Range("F14").Value = 1252
$applied++
$cur = $ws1.Range("F15").Value()
if ($cur -ne 4079) { $mismatches++ }
$ws1.Range("F15").Value = 4094
$applied++
$cur = $ws1.Range("F16").Value()
if ($cur -ne 957) { $mismatches++ }
$ws1.Range("F16").Value = 961
$applied++
$cur = $ws1.Range("F18").Value()
if ($cur -ne 1618) { $mismatches++ }
$ws1.Range("F18").Value = 1620
$applied++
$cur = $ws1.Range("F19").Value()
if ($cur -ne 2567) { $mismatches++ }
$ws1.Range("F19").Value = 2569
$applied++
$cur = $ws1.Range("F20").Value()
if ($cur -ne 18) { $mismatches++ }
$ws1.Range("F20").Value = 20
$applied++
$cur = $ws1.Range("F21").Value()
if ($cur -ne 95) { $mismatches++ }
$ws1.Range("F21").Value = 96
$applied++
$cur = $ws1.Range("F22").Value()
if ($cur -ne 136) { $mismatches++ }
$ws1.Range("F22").Value = 138
$applied++
$cur = $ws1.Range("F23").Value()
if ($cur -ne 931) { $mismatches++ }
$ws1.Range("F23").Value = 934
$applied++
$cur = $ws1.Range("F24").Value()
if ($cur -ne 274) { $mismatches++ }
$ws1.Range("F24").Value = 275
$applied++
$cur = $ws1.Range("F26").Value()
if ($cur -ne 69) { $mismatches++ }
$ws1.Range("F26").Value = 70
$applied++
$cur = $ws1.Range("F27").Value()
if ($cur -ne 193) { $mismatches++ }
$ws1.Range("F27").Value = 197
$applied++
$cur = $ws1.Range("F29").Value()
if ($cur -ne 319) { $mismatches++ }
$ws1.Range("F29").Value = 325
$applied++
$cur = $ws1.Range("F30").Value()
if ($cur -ne 0) { $mismatches++ }
$ws1.Range("F30").Value = 5
$applied++
$cur = $ws1.Range("F31").Value()
if ($cur -ne 96) { $mismatches++ }
$ws1.Range("F31").Value = 100
$applied++
$cur = $ws1.Range("F33").Value()
if ($cur -ne 177) { $mismatches++ }
$ws1.Range("F33").Value = 179
$applied++
$cur = $ws1.Range("F34").Value()
if ($cur -ne 1567) { $mismatches++ }
$ws1.Range("F34").Value = 1575
$applied++
$cur = $ws1.Range("F35").Value()
if ($cur -ne 2108) { $mismatches++ }
$ws1.Range("F35").Value = 2115
$applied++
$cur = $ws1.Range("F36").Value()
if ($cur -ne 990) { $mismatches++ }
$ws1.Range("F36").Value = 991
$applied++
$cur = $ws1.Range("F39").Value()
if ($cur -ne 572) { $mismatches++ }
$ws1.Range("F39").Value = 574
$applied++
$cur = $ws1.Range("F40").Value()
if ($cur -ne 222) { $mismatches++ }
$ws1.Range("F40").Value = 226
$applied++
$cur = $ws1.Range("I40").Value()
if ($cur -ne "//i1.hdslb.com/bfs/openplatform/202404/4rF9ZrcA1712820950457.jpeg") { $mismatches++ }
$ws1.Range("I40").Value = "//i1.hdslb.com/bfs/openplatform/202404/Na7jHnDL1713774453606.jpeg"
$applied++
$cur = $ws1.Range("F43").Value()
if ($cur -ne 367) { $mismatches++ }
$ws1.Range("F43").Value = 371
$applied++
$cur = $ws1.Range("F44").Value()
if ($cur -ne 255) { $mismatches++ }
$ws1.Range("F44").Value = 259
$applied++
$cur = $ws1.Range("F45").Value()
if ($cur -ne 194) { $mismatches++ }
$ws1.Range("F45").Value = 195
$applied++
$cur = $ws1.Range("F46").Value()
if ($cur -ne 114) { $mismatches++ }
$ws1.Range("F46").Value = 115
$applied++

$ws2 = $wb.Worksheets.Item("演出")
$cur = $ws2.Range("F7").Value()
if ($cur -ne 27) { $mismatches++ }
$ws2.Range("F7").Value = 28
$applied++
$cur = $ws2.Range("F10").Value()
if ($cur -ne 139) { $mismatches++ }
$ws2.Range("F10").Value = 141
$applied++

$ws3 = $wb.Worksheets.Item("本地生活")
$cur = $ws3.Range("F2").Value()
if ($cur -ne 680) { $mismatches++ }
$ws3.Range("F2").Value = 684
$applied++

$ws4 = $wb.Worksheets.Item("全部类型")
$cur = $ws4.Range("F2").Value()
if ($cur -ne 680) { $mismatches++ }
$ws4.Range("F2").Value = 684
$applied++
$cur = $ws4.Range("F7").Value()
if ($cur -ne 556) { $mismatches++ }
$ws4.Range("F7").Value = 557
$applied++
$cur = $ws4.Range("F8").Value()
if ($cur -ne 18) { $mismatches++ }
$ws4.Range("F8").Value = 20
$applied++
$cur = $ws4.Range("F9").Value()
if ($cur -ne 1907) { $mismatches++ }
$ws4.Range("F9").Value = 1910
$applied++
$cur = $ws4.Range("F10").Value()
if ($cur -ne 5051) { $mismatches++ }
$ws4.Range("F10").Value = 5068
$applied++
$cur = $ws4.Range("F11").Value()
if ($cur -ne 1412) { $mismatches++ }
$ws4.Range("F11").Value = 1419
$applied++
$cur = $ws4.Range("F12").Value()
if ($cur -ne 148) { $mismatches++ }
$ws4.Range("F12").Value = 149
$applied++
$cur = $ws4.Range("F14").Value()
if ($cur -ne 3005) { $mismatches++ }
$ws4.Range("F14").Value = 3010
$applied++
$cur = $ws4.Range("F15").Value()
if ($cur -ne 29) { $mismatches++ }
$ws4.Range("F15").Value = 31
$applied++
$cur = $ws4.Range("F16").Value()
if ($cur -ne 1250) { $mismatches++ }
$ws4.Range("F16").Value = 1252
$applied++
$cur = $ws4.Range("F17").Value()
if ($cur -ne 4079) { $mismatches++ }
$ws4.Range("F17").Value = 4095
$applied++
$cur = $ws4.Range("F18").Value()
if ($cur -ne 957) { $mismatches++ }
$ws4.Range("F18").Value = 961
$applied++
$cur = $ws4.Range("F20").Value()
if ($cur -ne 1618) { $mismatches++ }
$ws4.Range("F20").Value = 1620
$applied++
$cur = $ws4.Range("F22").Value()
if ($cur -ne 2567) { $mismatches++ }
$ws4.Range("F22").Value = 2569
$applied++
$cur = $ws4.Range("F23").Value()
if ($cur -ne 27) { $mismatches++ }
$ws4.Range("F23").Value = 28
$applied++
$cur = $ws4.Range("F24").Value()
if ($cur -ne 18) { $mismatches++ }
$ws4.Range("F24").Value = 20
$applied++
$cur = $ws4.Range("F27").Value()
if ($cur -ne 136) { $mismatches++ }
$ws4.Range("F27").Value = 138
$applied++
$cur = $ws4.Range("F28").Value()
if ($cur -ne 139) { $mismatches++ }
$ws4.Range("F28").Value = 141
$applied++
$cur = $ws4.Range("F29").Value()
if ($cur -ne 931) { $mismatches++ }
$ws4.Range("F29").Value = 934
$applied++
$cur = $ws4.Range("F30").Value()
if ($cur -ne 274) { $mismatches++ }
$ws4.Range("F30").Value = 275
$applied++
$cur = $ws4.Range("F31").Value()
if ($cur -ne 69) { $mismatches++ }
$ws4.Range("F31").Value = 70
$applied++
$cur = $ws4.Range("F32").Value()
if ($cur -ne 193) { $mismatches++ }
$ws4.Range("F32").Value = 197
$applied++
$cur = $ws4.Range("F35").Value()
if ($cur -ne 319) { $mismatches++ }
$ws4.Range("F35").Value = 325
$applied++
$cur = $ws4.Range("F36").Value()
if ($cur -ne 1567) { $mismatches++ }
$ws4.Range("F36").Value = 1575
$applied++
$cur = $ws4.Range("F37").Value()
if ($cur -ne 2108) { $mismatches++ }
$ws4.Range("F37").Value = 2115
$applied++
$cur = $ws4.Range("F39").Value()
if ($cur -ne 990) { $mismatches++ }
$ws4.Range("F39").Value = 991
$applied++
$cur = $ws4.Range("F43").Value()
if ($cur -ne 572) { $mismatches++ }
$ws4.Range("F43").Value = 574
$applied++
$cur = $ws4.Range("F44").Value()
if ($cur -ne 222) { $mismatches++ }
$ws4.Range("F44").Value = 226
$applied++
$cur = $ws4.Range("I44").Value()
if ($cur -ne "//i1.hdslb.com/bfs/openplatform/202404/4rF9ZrcA1712820950457.jpeg") { $mismatches++ }
$ws4.Range("I44").Value = "//i1.hdslb.com/bfs/openplatform/202404/Na7jHnDL1713774453606.jpeg"
$applied++
$cur = $ws4.Range("F46").Value()
if ($cur -ne 367) { $mismatches++ }
$ws4.Range("F46").Value = 371
$applied++
$cur = $ws4.Range("F47").Value()
if ($cur -ne 255) { $mismatches++ }
$ws4.Range("F47").Value = 259
$applied++
$cur = $ws4.Range("F48").Value()
if ($cur -ne 194) { $mismatches++ }
$ws4.Range("F48").Value = 195
$applied++
$cur = $ws4.Range("F49").Value()
if ($cur -ne 114) { $mismatches++ }
$ws4.Range("F49").Value = 115
$applied++

Write-Host ("applied=" + $applied + " mismatches=" + $mismatches)
